$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-17 14:41:50"
$wsZh.Range("G2").Value = "2016-01-17 14:43:19"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-17 14:42:00"
$wsDe.Range("G2").Value = "2016-01-17 14:43:37"
